$d = $word.ActiveDocument

$pairs = @(
    @("26×12=", "43×90="),
    @("77×53=", "50×95="),
    @("67×84=", "90×28="),
    @("14×69=", "59×63="),
    @("79×75=", "80×42="),
    @("25×78=", "71×54="),
    @("22×73=", "48×96="),
    @("80×86=", "37×41="),
    @("64×37=", "31×87="),
    @("23×67=", "56×54="),
    @("49×78=", "64×31="),
    @("76×95=", "24×56="),
    @("19×40=", "89×76="),
    @("29×74=", "84×29="),
    @("41×28=", "22×36="),
    @("16×81=", "66×66="),
    @("97×74=", "48×36="),
    @("62×15=", "28×35="),
    @("93×32=", "42×33="),
    @("32×56=", "95×38="),
    @("65×37=", "19×86="),
    @("27×23=", "18×35="),
    @("93×47=", "71×88="),
    @("38×99=", "53×55="),
    @("51×43=", "81×19=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
